# Update tag column (P) for rows where:
#   - load_pretrained_weights (column H) is TRUE
#   - tag (column P) is exactly "DescEmb-RNN_Scr"
# to "DescEmb-RNN_Scr-MLM"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $loadPretrained = $ws.Cells.Item($r, 8).Value2   # column H
    $tag = $ws.Cells.Item($r, 16).Value2             # column P

    if ($tag -eq "DescEmb-RNN_Scr" -and $loadPretrained -eq $true) {
        $ws.Cells.Item($r, 16).Value = "DescEmb-RNN_Scr-MLM"
    }
}
